$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8300.1
$ws.Range("I40").Value = 7571.5713
$ws.Range("K40").Value = 7571.5713
$ws.Range("M40").Value = -7396.5713

$ws.Range("H62").Value = 4809638
$ws.Range("I62").Value = 6251430
$ws.Range("K62").Value = 6251430
$ws.Range("M62").Value = -6250806

$ws.Range("H65").Value = 4809638
$ws.Range("I65").Value = 6251430
$ws.Range("K65").Value = 31257150
$ws.Range("M65").Value = -31254030

$ws.Range("H74").Value = 9662.823
$ws.Range("J74").Value = 20915
$ws.Range("L74").Value = 20915
$ws.Range("N74").Value = -22787

$ws.Range("H77").Value = 9662.823
$ws.Range("J77").Value = 20915
$ws.Range("L77").Value = 104575
$ws.Range("N77").Value = -113935

$ws.Range("H127").Value = 9796.846
$ws.Range("I127").Value = 19143.166
$ws.Range("J127").Value = 1785.7142
$ws.Range("K127").Value = 57429.49800000001
$ws.Range("L127").Value = 5357.142599999999
$ws.Range("M127").Value = -52469.49800000001
$ws.Range("N127").Value = -15277.1426

$ws.Range("H129").Value = 11130.167
$ws.Range("I129").Value = 12898.223
$ws.Range("K129").Value = 38694.669
$ws.Range("M129").Value = -33694.669

$ws.Range("H137").Value = 6072.607
$ws.Range("I137").Value = 3804.9167
$ws.Range("K137").Value = 11414.7501
$ws.Range("M137").Value = -8864.750100000001

$ws.Range("H138").Value = 6732.323
$ws.Range("I138").Value = 5832.2
$ws.Range("K138").Value = 17496.6
$ws.Range("M138").Value = -12356.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 3142.4443
$ws.Range("I50").Value = 4011.75
$ws.Range("K50").Value = 4011.75
$ws.Range("M50").Value = -3297.75

$ws.Range("H61").Value = 4642.6816
$ws.Range("I61").Value = 2856.8462
$ws.Range("K61").Value = 2856.8462
$ws.Range("M61").Value = -2644.8462

$ws.Range("H97").Value = 3962.5715
$ws.Range("I97").Value = 3962.5715
$ws.Range("K97").Value = 3962.5715
$ws.Range("M97").Value = -3466.5715

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774

$ws.Range("H122").Value = 5325
$ws.Range("I122").Value = 2955
$ws.Range("K122").Value = 8865
$ws.Range("M122").Value = -6415

$ws.Range("H132").Value = 3179.6
$ws.Range("I132").Value = 1529.1471
$ws.Range("K132").Value = 4587.4413
$ws.Range("M132").Value = -2057.4413

$ws.Range("H136").Value = 4642.6816
$ws.Range("I136").Value = 2856.8462
$ws.Range("K136").Value = 8570.5386
$ws.Range("M136").Value = -6020.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 59998.6
$ws.Range("J137").Value = 59998.6
$ws.Range("L137").Value = 59998.6
$ws.Range("N137").Value = -70198.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 4207.231
$ws.Range("I15").Value = 1447.6666
$ws.Range("J15").Value = 6572.5713
$ws.Range("K15").Value = 1447.6666
$ws.Range("L15").Value = 6572.5713
$ws.Range("M15").Value = -1277.6666
$ws.Range("N15").Value = -6912.5713

$ws.Range("H31").Value = 361487.44
$ws.Range("I31").Value = 557825.9399999999
$ws.Range("K31").Value = 557825.9399999999
$ws.Range("M31").Value = -557530.9399999999

$ws.Range("H34").Value = 361487.44
$ws.Range("I34").Value = 557825.9399999999
$ws.Range("K34").Value = 557825.9399999999
$ws.Range("M34").Value = -557623.9399999999

$ws.Range("H100").Value = 56331.668
$ws.Range("J100").Value = 56331.668
$ws.Range("L100").Value = 56331.668
$ws.Range("N100").Value = -58495.668

$ws.Range("H107").Value = 849.4545000000001
$ws.Range("I107").Value = 849.4545000000001
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 849.4545000000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1070.5455
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 2706.7932
$ws.Range("I122").Value = 2127.5
$ws.Range("K122").Value = 6382.5
$ws.Range("M122").Value = -3932.5

$ws.Range("H134").Value = 299048.6
$ws.Range("I134").Value = 3647.611
$ws.Range("K134").Value = 10942.833
$ws.Range("M134").Value = -8407.832999999999

$ws.Range("H141").Value = 101247.586
$ws.Range("I141").Value = 45000
$ws.Range("J141").Value = 106361
$ws.Range("K141").Value = 45000
$ws.Range("L141").Value = 106361
$ws.Range("M141").Value = -39820
$ws.Range("N141").Value = -116721

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 81769
$ws.Range("J5").Value = 599
$ws.Range("L5").Value = 1797
$ws.Range("N5").Value = -2021

$ws.Range("H55").Value = 9849.5625
$ws.Range("I55").Value = 1242.5714
$ws.Range("K55").Value = 3727.7142
$ws.Range("M55").Value = -3550.7142

$ws.Range("H98").Value = 1819.7778
$ws.Range("I98").Value = 2321.75
$ws.Range("K98").Value = 6965.25
$ws.Range("M98").Value = -5467.25

$ws.Range("H103").Value = 265.6
$ws.Range("I103").Value = 208.5
$ws.Range("J103").Value = 303.66666
$ws.Range("K103").Value = 625.5
$ws.Range("L103").Value = 910.9999799999999
$ws.Range("M103").Value = 253.5
$ws.Range("N103").Value = -2668.99998

$ws.Range("H131").Value = 100040.14
$ws.Range("J131").Value = 57781.26
$ws.Range("L131").Value = 173343.78
$ws.Range("N131").Value = -183423.78

$ws.Range("H135").Value = 81769
$ws.Range("J135").Value = 599
$ws.Range("L135").Value = 5391
$ws.Range("N135").Value = -10461

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1665.4642
$ws.Range("I97").Value = 1575.3636
$ws.Range("K97").Value = 1575.3636
$ws.Range("M97").Value = -1079.3636

$ws.Range("H102").Value = 1433.5625
$ws.Range("I102").Value = 1040.2307
$ws.Range("J102").Value = 3138
$ws.Range("K102").Value = 1040.2307
$ws.Range("L102").Value = 3138
$ws.Range("M102").Value = 581.7692999999999
$ws.Range("N102").Value = -6382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 563509.3
$ws.Range("I7").Value = 10813.857
$ws.Range("J7").Value = 915224.6
$ws.Range("K7").Value = 10813.857
$ws.Range("L7").Value = 915224.6
$ws.Range("M7").Value = -10701.857
$ws.Range("N7").Value = -915448.6

$ws.Range("H126").Value = 563509.3
$ws.Range("I126").Value = 10813.857
$ws.Range("J126").Value = 915224.6
$ws.Range("K126").Value = 32441.571
$ws.Range("L126").Value = 2745673.8
$ws.Range("M126").Value = -29971.571
$ws.Range("N126").Value = -2750613.8

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H136").Value = 4071.9333
$ws.Range("I136").Value = 2827.9048
$ws.Range("J136").Value = 6974.6665
$ws.Range("K136").Value = 8483.714399999999
$ws.Range("L136").Value = 20923.9995
$ws.Range("M136").Value = -5933.714399999999
$ws.Range("N136").Value = -26023.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1293.8
$ws.Range("J113").Value = 1978.3334
$ws.Range("L113").Value = 5935.0002
$ws.Range("N113").Value = -10275.0002

$ws.Range("H136").Value = 60471.465
$ws.Range("I136").Value = 12823.743
$ws.Range("K136").Value = 38471.229
$ws.Range("M136").Value = -35921.229
